$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.497.45'
$ws.Range("D3").Value = '1.909.53'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''243.93'
$ws.Range("D7").Value = '''0.4834'
$ws.Range("D8").Value = '''0.2883'
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("D10").Value = '''111.49'
$ws.Range("E10").Value = '  +6.16%  '
$ws.Range("D11").Value = '''19.32'
$ws.Range("E11").Value = '  +5.22%  '
$ws.Range("D12").Value = '1.915.10'
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").Value = '''0.07563'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '''5.386'
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").Value = '''0.6679'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").Value = '''292.30'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").Value = '30.499.82'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = '''12.98'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").Value = '2.163.14'
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '''6.388'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").Value = '''9.442'
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("D26").Value = '''165.29'
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("D27").Value = '''20.21'
$ws.Range("E27").Value = '  -4.14%  '
$ws.Range("D28").Value = '''2.065'
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("D29").Value = '''0.1065'
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").Value = '''1.429'
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("D31").Value = '''4.123'
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("D32").Value = '''4.047'
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").Value = '''0.04969'
$ws.Range("E33").Value = '  -1.12%  '
$ws.Range("D34").Value = '''0.7345'
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.02028'
$ws.Range("E37").Value = '  -2.03%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '''2.708'
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("D39").Value = '''2.681'
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''2.008'
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '''109.26'
$ws.Range("E41").Value = '  -1.62%  '
$ws.Range("D42").Value = '''0.4418'
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").Value = '''0.8611'
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("D44").Value = '''5.776'
$ws.Range("E44").Value = '  -1.55%  '
$ws.Range("D45").Value = '''1.001'
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '''68.98'
$ws.Range("E46").Value = '  +2.20%  '
$ws.Range("D47").Value = '''7.176'
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").Value = '''48.04'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''9.209'
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("D50").Value = '''0.1225'
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").Value = '''0.2514'
$ws.Range("E51").Value = '  +1.40%  '
